$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.959.52"
$ws.Range("E2").Value = "  -3.19%  "
$ws.Range("D3").Value = "2.288.05"
$ws.Range("E3").Value = "  -3.69%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.79"
$ws.Range("E5").Value = "  -4.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.16"
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "2.286.45"
$ws.Range("E10").Value = "  -5.67%  "
$ws.Range("E11").Value = "  -3.81%  "
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.53"
$ws.Range("E14").Value = "  -3.72%  "
$ws.Range("D15").Value = "2.696.71"
$ws.Range("E15").Value = "  -3.60%  "
$ws.Range("D16").Value = "57.921.73"
$ws.Range("E16").Value = "  -3.18%  "
$ws.Range("E17").Value = "  -4.52%  "
$ws.Range("D18").Value = "2.281.31"
$ws.Range("E18").Value = "  -4.05%  "
$ws.Range("E19").Value = "  -5.71%  "
$ws.Range("E20").Value = "  -5.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.90"
$ws.Range("E21").Value = "  -2.89%  "
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.37"
$ws.Range("E24").Value = "  -2.77%  "
$ws.Range("E25").Value = "  -3.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -5.45%  "
$ws.Range("E28").Value = "  -6.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.83"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.70"
$ws.Range("E30").Value = "  -5.93%  "
$ws.Range("D31").Value = "0.0₃0719"
$ws.Range("E31").Value = "  -5.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.76"
$ws.Range("E32").Value = "  -5.23%  "
$ws.Range("E33").Value = "  -7.56%  "
$ws.Range("E34").Value = "  -5.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.73"
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  -7.33%  "
$ws.Range("E39").Value = "  -5.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.13"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("E41").Value = "  -6.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "141.72"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "287.92"
$ws.Range("E43").Value = "  -9.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.43"
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0949"
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("E47").Value = "  -2.69%  "
$ws.Range("E48").Value = "  -8.11%  "
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("E51").Value = "  -0.75%  "
